# Add a new "student_type" column (K) to the roster sheet.
# Header in K1, and a value of "New Student" or "Transferee" for every
# data row (2-120). Five specific students (rows 43, 58, 87, 104, 120)
# are marked "Transferee"; all the rest are "New Student".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("K1").Value = "student_type"

# Rows that should read "Transferee" instead of "New Student"
$transfereeRows = @(43, 58, 87, 104, 120)

for ($r = 2; $r -le 120; $r++) {
    if ($transfereeRows -contains $r) {
        $ws.Cells.Item($r, 11).Value = "Transferee"
    } else {
        $ws.Cells.Item($r, 11).Value = "New Student"
    }
}

# Size the new column to fit its contents, like Excel does automatically
# after typing into a fresh column.
$ws.Columns.Item(11).AutoFit()

# Update the view: active selection matches the author's final
# on-screen state.
$ws.Range("P50").Select()
